$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content (keeps row 1 header formatting, resets shared-string table)
$ws.Cells.ClearContents()

# --- Header row (rebuild shared strings 0..19 in original order) ---
$ws.Range("A1").Value = 'Sending cluster'
$ws.Range("B1").Value = 'Ligand symbol'
$ws.Range("C1").Value = 'Receptor symbol'
$ws.Range("D1").Value = 'Target cluster'
$ws.Range("E1").Value = 'Ligand-expressing cells'
$ws.Range("F1").Value = 'Ligand detection rate'
$ws.Range("G1").Value = 'Ligand average expression value'
$ws.Range("H1").Value = 'Ligand total expression value'
$ws.Range("I1").Value = 'Ligand derived specificity of average expression value'
$ws.Range("J1").Value = 'Ligand derived specificity of total expression value'
$ws.Range("K1").Value = 'Receptor-expressing cells'
$ws.Range("L1").Value = 'Receptor detection rate'
$ws.Range("M1").Value = 'Receptor average expression value'
$ws.Range("N1").Value = 'Receptor total expression value'
$ws.Range("O1").Value = 'Receptor derived specificity of average expression value'
$ws.Range("P1").Value = 'Receptor derived specificity of total expression value'
$ws.Range("Q1").Value = 'Edge average expression weight'
$ws.Range("R1").Value = 'Edge total expression weight'
$ws.Range("S1").Value = 'Edge average expression derived specificity'
$ws.Range("T1").Value = 'Edge total expression derived specificity'

# --- Column A (Sending cluster) rows 2-13: establishes ECs, FAPs, Resolving-Mac order ---
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"
$ws.Range("A5").Value = "ECs"
$ws.Range("A6").Value = "FAPs"
$ws.Range("A7").Value = "FAPs"
$ws.Range("A8").Value = "FAPs"
$ws.Range("A9").Value = "FAPs"
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("A13").Value = "Resolving-Mac"

# --- Column B (Ligand symbol) rows 2-13: establishes Lta ---
$ws.Range("B2").Value = "Lta"
$ws.Range("B3").Value = "Lta"
$ws.Range("B4").Value = "Lta"
$ws.Range("B5").Value = "Lta"
$ws.Range("B6").Value = "Lta"
$ws.Range("B7").Value = "Lta"
$ws.Range("B8").Value = "Lta"
$ws.Range("B9").Value = "Lta"
$ws.Range("B10").Value = "Lta"
$ws.Range("B11").Value = "Lta"
$ws.Range("B12").Value = "Lta"
$ws.Range("B13").Value = "Lta"

# --- Column C (Receptor symbol) rows 2-13: establishes Tnfrsf1b ---
$ws.Range("C2").Value = "Tnfrsf1b"
$ws.Range("C3").Value = "Tnfrsf1b"
$ws.Range("C4").Value = "Tnfrsf1b"
$ws.Range("C5").Value = "Tnfrsf1b"
$ws.Range("C6").Value = "Tnfrsf1b"
$ws.Range("C7").Value = "Tnfrsf1b"
$ws.Range("C8").Value = "Tnfrsf1b"
$ws.Range("C9").Value = "Tnfrsf1b"
$ws.Range("C10").Value = "Tnfrsf1b"
$ws.Range("C11").Value = "Tnfrsf1b"
$ws.Range("C12").Value = "Tnfrsf1b"
$ws.Range("C13").Value = "Tnfrsf1b"

# --- Column D (Target cluster) rows 2-13: establishes MuSCs (last new string) ---
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("D6").Value = "ECs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("D10").Value = "ECs"
$ws.Range("D11").Value = "FAPs"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("D13").Value = "Resolving-Mac"

# --- Numeric columns E..T for each row ---
# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.103879
$ws.Range("H2").Value = 0.311637
$ws.Range("I2").Value = 0.1899090787212519
$ws.Range("J2").Value = 0.1899090787212519
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.57753066666667
$ws.Range("N2").Value = 37.732592
$ws.Range("O2").Value = 0.1317204310459389
$ws.Range("P2").Value = 0.1317204310459389
$ws.Range("Q2").Value = 1.306541308122667
$ws.Range("R2").Value = 11.758871773104
$ws.Range("S2").Value = 0.02501490570870045
$ws.Range("T2").Value = 0.02501490570870045

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.103879
$ws.Range("H3").Value = 0.311637
$ws.Range("I3").Value = 0.1899090787212519
$ws.Range("J3").Value = 0.1899090787212519
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.940628
$ws.Range("N3").Value = 32.821884
$ws.Range("O3").Value = 0.1145776761962127
$ws.Range("P3").Value = 0.1145776761962127
$ws.Range("Q3").Value = 1.136501496012
$ws.Range("R3").Value = 10.228513464108
$ws.Range("S3").Value = 0.02175934092844467
$ws.Range("T3").Value = 0.02175934092844467

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.103879
$ws.Range("H4").Value = 0.311637
$ws.Range("I4").Value = 0.1899090787212519
$ws.Range("J4").Value = 0.1899090787212519
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.520244333333333
$ws.Range("N4").Value = 7.560733
$ws.Range("O4").Value = 0.02639370785296846
$ws.Range("P4").Value = 0.02639370785296846
$ws.Range("Q4").Value = 0.2618004611023333
$ws.Range("R4").Value = 2.356204149921
$ws.Range("S4").Value = 0.005012404742395112
$ws.Range("T4").Value = 0.005012404742395112

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.103879
$ws.Range("H5").Value = 0.311637
$ws.Range("I5").Value = 0.1899090787212519
$ws.Range("J5").Value = 0.1899090787212519
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 69.44815566666666
$ws.Range("N5").Value = 208.344467
$ws.Range("O5").Value = 0.7273081849048799
$ws.Range("P5").Value = 0.7273081849048799
$ws.Range("Q5").Value = 7.214204962497666
$ws.Range("R5").Value = 64.927844662479
$ws.Range("S5").Value = 0.1381224273417117
$ws.Range("T5").Value = 0.1381224273417117

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.257314
$ws.Range("H6").Value = 0.771942
$ws.Range("I6").Value = 0.4704152396738534
$ws.Range("J6").Value = 0.4704152396738534
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.57753066666667
$ws.Range("N6").Value = 37.732592
$ws.Range("O6").Value = 0.1317204310459389
$ws.Range("P6").Value = 0.1317204310459389
$ws.Range("Q6").Value = 3.236374725962667
$ws.Range("R6").Value = 29.127372533664
$ws.Range("S6").Value = 0.06196329814041862
$ws.Range("T6").Value = 0.06196329814041862

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.257314
$ws.Range("H7").Value = 0.771942
$ws.Range("I7").Value = 0.4704152396738534
$ws.Range("J7").Value = 0.4704152396738534
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.940628
$ws.Range("N7").Value = 32.821884
$ws.Range("O7").Value = 0.1145776761962127
$ws.Range("P7").Value = 0.1145776761962127
$ws.Range("Q7").Value = 2.815176753192
$ws.Range("R7").Value = 25.336590778728
$ws.Range("S7").Value = 0.05389908500911456
$ws.Range("T7").Value = 0.05389908500911456

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.257314
$ws.Range("H8").Value = 0.771942
$ws.Range("I8").Value = 0.4704152396738534
$ws.Range("J8").Value = 0.4704152396738534
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.520244333333333
$ws.Range("N8").Value = 7.560733
$ws.Range("O8").Value = 0.02639370785296846
$ws.Range("P8").Value = 0.02639370785296846
$ws.Range("Q8").Value = 0.6484941503873333
$ws.Range("R8").Value = 5.836447353486
$ws.Range("S8").Value = 0.01241600240553582
$ws.Range("T8").Value = 0.01241600240553582

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.257314
$ws.Range("H9").Value = 0.771942
$ws.Range("I9").Value = 0.4704152396738534
$ws.Range("J9").Value = 0.4704152396738534
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 69.44815566666666
$ws.Range("N9").Value = 208.344467
$ws.Range("O9").Value = 0.7273081849048799
$ws.Range("P9").Value = 0.7273081849048799
$ws.Range("Q9").Value = 17.86998272721267
$ws.Range("R9").Value = 160.829844544914
$ws.Range("S9").Value = 0.3421368541187844
$ws.Range("T9").Value = 0.3421368541187844

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1858003333333333
$ws.Range("H10").Value = 0.557401
$ws.Range("I10").Value = 0.3396756816048946
$ws.Range("J10").Value = 0.3396756816048946
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.57753066666667
$ws.Range("N10").Value = 37.732592
$ws.Range("O10").Value = 0.1317204310459389
$ws.Range("P10").Value = 0.1317204310459389
$ws.Range("Q10").Value = 2.336909390376889
$ws.Range("R10").Value = 21.032184513392
$ws.Range("S10").Value = 0.04474222719681982
$ws.Range("T10").Value = 0.04474222719681982

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1858003333333333
$ws.Range("H11").Value = 0.557401
$ws.Range("I11").Value = 0.3396756816048946
$ws.Range("J11").Value = 0.3396756816048946
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 10.940628
$ws.Range("N11").Value = 32.821884
$ws.Range("O11").Value = 0.1145776761962127
$ws.Range("P11").Value = 0.1145776761962127
$ws.Range("Q11").Value = 2.032772329276
$ws.Range("R11").Value = 18.294950963484
$ws.Range("S11").Value = 0.03891925025865346
$ws.Range("T11").Value = 0.03891925025865346

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1858003333333333
$ws.Range("H12").Value = 0.557401
$ws.Range("I12").Value = 0.3396756816048946
$ws.Range("J12").Value = 0.3396756816048946
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.520244333333333
$ws.Range("N12").Value = 7.560733
$ws.Range("O12").Value = 0.02639370785296846
$ws.Range("P12").Value = 0.02639370785296846
$ws.Range("Q12").Value = 0.4682622372147778
$ws.Range("R12").Value = 4.214360134933
$ws.Range("S12").Value = 0.008965300705037522
$ws.Range("T12").Value = 0.008965300705037522

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1858003333333333
$ws.Range("H13").Value = 0.557401
$ws.Range("I13").Value = 0.3396756816048946
$ws.Range("J13").Value = 0.3396756816048946
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 69.44815566666666
$ws.Range("N13").Value = 208.344467
$ws.Range("O13").Value = 0.7273081849048799
$ws.Range("P13").Value = 0.7273081849048799
$ws.Range("Q13").Value = 12.90349047225189
$ws.Range("R13").Value = 116.131414250267
$ws.Range("S13").Value = 0.2470489034443838
$ws.Range("T13").Value = 0.2470489034443838

